$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing "Sharp Optical Dust Sensor" row (row 11): new supplier/cost/link ---
$ws.Range("D11").Value = 9.7
$ws.Range("E11").Value = "AliExpress (HS Electronics)"
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = "http://www.aliexpress.com/item/Free-Shipping-PM2-5-GP2Y1010AU0F-SENSOR-AIR-QUALITY-DUST-GP2Y10-Compact-Optical-Dust-Sensor-Smoke-Particle/32314774144.html?spm=2114.01010208.3.47.T9wGdx&ws_ab_test=searchweb201556_0,searchweb201602_3_10017_10005_10006_10034_10021_507_10022_10020_10018_10019,searchweb201603_9&btsid=0f96a037-dd7a-4cb8-8149-8bec116ef816"

# --- New row 12: Plantower PMS5003 Laser Dust Sensor ---
$ws.Range("C12").Value = "Plantower PMS5003 Laser Dust Sensor"
$ws.Range("D12").Value = 47
$ws.Range("E12").Value = "AliExpress (HS Electronics)"
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = "http://www.aliexpress.com/item/PLANTOWER-PM2-5-SENSOR-laser-dust-sensor-G5-PMS5003-High-precision-laser-dust-concentration-sensor-digital/32618735056.html?spm=2114.01010208.3.10.wbvLWe&ws_ab_test=searchweb201556_0,searchweb201602_3_10017_10005_10006_10034_10021_507_10022_10020_10018_10019,searchweb201603_9&btsid=8b09eea3-1bb7-4347-b788-b4bbf9f0939b"

# --- New row 13: Nova SDS011 Laser Dust Sensor ---
$ws.Range("C13").Value = "Nova SDS011 Laser Dust Sensor"
$ws.Range("D13").Value = 36.4
$ws.Range("E13").Value = "AliExpress (HS Electronics)"
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = "http://www.aliexpress.com/item/nova-PM-sensor-SDS011-High-precision-laser-pm2-5-air-quality-detection-sensor-module-Super-dust/32617788139.html?spm=2114.01010208.3.20.gbKqjW&ws_ab_test=searchweb201556_0,searchweb201602_3_10017_10005_10006_10034_10021_507_10022_10020_10018_10019,searchweb201603_9&btsid=0144ffea-7eee-4de0-8be5-b055fc86186a"

# --- Extend the "Order 1" total formula to cover the new rows 9:13 ---
$ws.Range("E5").Formula = "=SUM(D9:D13,F9:F13)"

# --- Bold the header row of Order 1's table ---
$ws.Range("C8:G8").Font.Bold = $true

# --- Column widths (approximate values; the runtime quantizes to 1/6 char units) ---
$ws.Columns.Item(3).ColumnWidth = 32.666666666666664
$ws.Columns.Item(5).ColumnWidth = 22.666666666666668
$ws.Columns.Item(6).ColumnWidth = 12.5
$ws.Columns.Item(7).ColumnWidth = 50.666666666666664
$ws.Columns.Item(10).ColumnWidth = 13.0
$ws.Columns.Item(13).ColumnWidth = 12.5

# --- View: zoom to 85% and move the active selection ---
$ws.Range("E18").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
